$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = 1.33
$ws.Range("M2").Value = 3.4

# Row 3
$ws.Range("G3").Value = 3.1
$ws.Range("H3").Value = 3.25
$ws.Range("I3").Value = 2.4
$ws.Range("N3").Value = 2.01
$ws.Range("O3").Value = 1.89
$ws.Range("T3").Value = 10
$ws.Range("Z3").Value = 10

# Row 4
$ws.Range("G4").Value = 7.5
$ws.Range("H4").Value = 5.5
$ws.Range("I4").Value = 1.33
$ws.Range("T4").Value = 23
$ws.Range("V4").Value = 21
$ws.Range("W4").Value = 81
$ws.Range("AD4").Value = 151
$ws.Range("AH4").Value = 9.5

# Row 5
$ws.Range("N5").Value = 3.2
$ws.Range("O5").Value = 1.36

# Row 7
$ws.Range("G7").Value = 1.62
$ws.Range("H7").Value = 4.1
$ws.Range("J7").Value = 1.03
$ws.Range("K7").Value = 15
$ws.Range("N7").Value = 1.7
$ws.Range("O7").Value = 2.1
$ws.Range("W7").Value = 13
$ws.Range("X7").Value = 13
$ws.Range("AI7").Value = 34

# Row 8
$ws.Range("G8").Value = 1.33
$ws.Range("H8").Value = 5.25
$ws.Range("J8").Value = 1.02
$ws.Range("K8").Value = 21
$ws.Range("L8").Value = 1.11
$ws.Range("M8").Value = 6.5
$ws.Range("N8").Value = 1.4
$ws.Range("O8").Value = 2.88
$ws.Range("P8").Value = 1.22
$ws.Range("Q8").Value = 4
$ws.Range("R8").Value = 1.67
$ws.Range("S8").Value = 2.1
$ws.Range("T8").Value = 10
$ws.Range("U8").Value = 8.5
$ws.Range("W8").Value = 10
$ws.Range("X8").Value = 10
$ws.Range("Y8").Value = 19
$ws.Range("Z8").Value = 21
$ws.Range("AA8").Value = 11
$ws.Range("AB8").Value = 17
$ws.Range("AC8").Value = 41
$ws.Range("AD8").Value = 151
$ws.Range("AE8").Value = 26
$ws.Range("AG8").Value = 23
$ws.Range("AJ8").Value = 41

# Row 12
$ws.Range("G12").Value = 1.72
$ws.Range("I12").Value = 5
$ws.Range("R12").Value = 2.22
$ws.Range("X12").Value = 17.5
$ws.Range("AG12").Value = 18
$ws.Range("AI12").Value = 70
$ws.Range("AJ12").Value = 90

# Row 17
$ws.Range("H17").Value = 3
$ws.Range("N17").Value = 2.35
$ws.Range("O17").Value = 1.57
$ws.Range("P17").Value = 1.5
$ws.Range("Q17").Value = 2.5
$ws.Range("T17").Value = 6.5

# Row 23
$ws.Range("H23").Value = 3.4

# Row 24
$ws.Range("G24").Value = 3.5
$ws.Range("H24").Value = 3.5
$ws.Range("I24").Value = 2.05

# Row 32
$ws.Range("N32").Value = 1.75
$ws.Range("O32").Value = 2.05

# Row 33
$ws.Range("L33").Value = 1.29
$ws.Range("M33").Value = 3.5
$ws.Range("N33").Value = 1.9
$ws.Range("O33").Value = 1.95
